$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row suffixes: "_old" -> "_FV2210", "_new" -> "_FV2304"
$headerRange = $ws.Range("A1:U1")
$headerRange.Replace("_old", "_FV2210") | Out-Null
$headerRange.Replace("_new", "_FV2304") | Out-Null

# 2. Turn the data range into an Excel Table ("ListObject") named Table1
$dataRange = $ws.Range("A1:U87")
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"

# 3. Freeze the header row (freeze pane below row 1)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
